# Weekly update: insert this week's Cebollín price record at the top of the
# data block (row 176) and push the existing rows 176..244 down by one,
# appending the former last row (244) as a brand-new row 245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 176
$lastRow  = 244
$newRow   = 245

$dCol = 4
$jCol = 10
$kCol = 11
$lCol = 12
$mCol = 13
$nCol = 14
$pCol = 16
$qCol = 17

# 1) Snapshot the columns that carry per-row data (D,J,K,L,M,N,P,Q) for every
#    row in the existing block before we start overwriting anything.
$snapD = @{}
$snapJ = @{}
$snapK = @{}
$snapL = @{}
$snapM = @{}
$snapN = @{}
$snapP = @{}
$snapQ = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, $dCol).Value2()
    $snapJ[$r] = $ws.Cells.Item($r, $jCol).Value()
    $snapK[$r] = $ws.Cells.Item($r, $kCol).Value()
    $snapL[$r] = $ws.Cells.Item($r, $lCol).Value()
    $snapM[$r] = $ws.Cells.Item($r, $mCol).Value()
    $snapN[$r] = $ws.Cells.Item($r, $nCol).Value()
    $snapP[$r] = $ws.Cells.Item($r, $pCol).Value()
    $snapQ[$r] = $ws.Cells.Item($r, $qCol).Value()
}

# 2) Shift rows 177..244 to hold what used to be in the row just above them
#    (176..243), going from the bottom up so we never read an already
#    overwritten cell.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, $dCol).Value = $snapD[$src]
    $ws.Cells.Item($r, $jCol).Value = $snapJ[$src]
    $ws.Cells.Item($r, $kCol).Value = $snapK[$src]
    $ws.Cells.Item($r, $lCol).Value = $snapL[$src]
    $ws.Cells.Item($r, $mCol).Value = $snapM[$src]
    $ws.Cells.Item($r, $nCol).Value = $snapN[$src]
    $ws.Cells.Item($r, $pCol).Value = $snapP[$src]
    $ws.Cells.Item($r, $qCol).Value = $snapQ[$src]
}

# 3) Row 176 becomes this week's new record.
$ws.Cells.Item($firstRow, $dCol).Value = 44837
$ws.Cells.Item($firstRow, $jCol).Value = 1200
$ws.Cells.Item($firstRow, $kCol).Value = 1400
$ws.Cells.Item($firstRow, $lCol).Value = 1600
$ws.Cells.Item($firstRow, $mCol).Value = 1500
$ws.Cells.Item($firstRow, $nCol).Value = "`$/paquete 6 unidades"
$ws.Cells.Item($firstRow, $pCol).Value = 250
$ws.Cells.Item($firstRow, $qCol).Value = 6

# 4) Append a brand-new row 245 holding the former row 244 data (the record
#    that fell off the bottom of the shift). The descriptive columns
#    (A,B,C,E,F,G,H,I,O,R) are identical across the whole block, so copy them
#    straight from row 244.
for ($c = 1; $c -le 18; $c++) {
    if ($c -eq $dCol -or $c -eq $jCol -or $c -eq $kCol -or $c -eq $lCol -or $c -eq $mCol -or $c -eq $nCol -or $c -eq $pCol -or $c -eq $qCol) {
        continue
    }
    $ws.Cells.Item($newRow, $c).Value = $ws.Cells.Item($lastRow, $c).Value()
}

$ws.Cells.Item($newRow, $dCol).Value = $snapD[$lastRow]
$ws.Cells.Item($newRow, $jCol).Value = $snapJ[$lastRow]
$ws.Cells.Item($newRow, $kCol).Value = $snapK[$lastRow]
$ws.Cells.Item($newRow, $lCol).Value = $snapL[$lastRow]
$ws.Cells.Item($newRow, $mCol).Value = $snapM[$lastRow]
$ws.Cells.Item($newRow, $nCol).Value = $snapN[$lastRow]
$ws.Cells.Item($newRow, $pCol).Value = $snapP[$lastRow]
$ws.Cells.Item($newRow, $qCol).Value = $snapQ[$lastRow]

# Match the date formatting used by the rest of column D.
$ws.Cells.Item($newRow, $dCol).NumberFormat = $ws.Cells.Item($lastRow, $dCol).NumberFormat
